$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump the generated "Date" timestamp -----------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- 2. Elements sheet: add a new mapping column (draft mapping) ------------
$ws = $wb.Worksheets.Item("Elements")

# Copy the formatting of the existing last "Mapping" column (AK) onto the
# new column (AL) so header/body styling (border, bold header, wrap, etc.)
# matches the rest of the mapping columns.
$ws.Range("AK1").Copy()
$ws.Range("AL1").PasteSpecial(-4122)
$ws.Range("AK2:AK6").Copy()
$ws.Range("AL2:AL6").PasteSpecial(-4122)

# Header for the new mapping column.
$ws.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR TelecomConfidentialityLevel"

# Only the Extension.value[x] row (row 6) has a mapped value so far.
$ws.Range("AL6").Value = "niveauConfidentialite"

# Widen the new column to fit its (long) content, mirroring the other
# bestFit mapping/description columns.
$ws.Columns.Item(38).ColumnWidth = 83.5
